$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.455.07'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.161.67'
$ws.Range("E3").Value = '  +3.19%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.82'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.42'
$ws.Range("E7").Value = '  +4.65%  '
$ws.Range("E9").Value = '  +2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0860'
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.04'
$ws.Range("E12").Value = '  +4.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.481.43'
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.24'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.815'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.165.67'
$ws.Range("E17").Value = '  +3.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.415.47'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0853'
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.67'
$ws.Range("E22").Value = '  +1.55%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +5.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.47'
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.52'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("E29").Value = '  +2.83%  '
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  +8.84%  '
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.62'
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.75'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("E35").Value = '  +8.74%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.41'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.58'
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '104.02'
$ws.Range("E40").Value = '  +2.55%  '
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.89'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.539.23'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.19'
$ws.Range("E44").Value = '  +4.37%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0926'
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  +5.76%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.81'
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.365.68'
$ws.Range("E50").Value = '  +3.27%  '
$ws.Range("E51").Value = '  +0.01%  '
